# profile onboarding script implementation
# Adds a new "Profile60" / OPQA-2096 test-case row (row 61) to the
# "Test Cases" sheet, mirroring the formatting of the row above it,
# and moves the active selection to C46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 61 with the same cell formatting as the template row above it
# (row 60), then fix up the Runmode cell (D61) to reuse the formatting
# already used for that column elsewhere in the sheet (B58).
$ws.Range("A60:E60").Copy() | Out-Null
$ws.Range("A61:E61").PasteSpecial(-4122) | Out-Null

$ws.Range("B58").Copy() | Out-Null
$ws.Range("D61").PasteSpecial(-4122) | Out-Null

# Populate the new test case.
$ws.Range("A61").Value = "Profile60"
$ws.Range("B61").Value = "OPQA-2096"
$ws.Range("C61").Value = "Verify that the system records the user as on-boarded if he exits the Neon on-boarding welcome modal without clicking on the ""Done"" button."
$ws.Range("D61").Value = "Y"

# Match the author's final selection in the sheet.
$ws.Range("C46").Select() | Out-Null
